# Applies the "update oa + doc gen" edit:
#  1. "Mr/Ms {client}" -> split into "{gender}" run + " {client}" run
#  2. "Dear Mr/Ms {client} , " -> "Dear Mr/Ms , "
#  3. merge " (", "{TVA}", "), as retainer regarding your case. " runs
#     into a single run with text " ({TVA}), as retainer regarding your case. "

$d = $word.ActiveDocument

# --- Change 2 (done first so its "Mr/Ms {client}" text stops matching
#     the Change-1 search below) ---------------------------------------
$d.Content.Find.Execute("Dear Mr/Ms {client} , ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Dear Mr/Ms , ", 2) | Out-Null

# --- Change 1 -------------------------------------------------------
# Replace the whole "Mr/Ms {client}" run text first.
$d.Content.Find.Execute("Mr/Ms {client}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{gender} {client}", 2) | Out-Null

# Re-find just the " {client}" tail and toggle a character property on
# and back off again; this forces the COM host to split it into its own
# run (with identical rPr to its neighbour) without altering formatting.
$tail = $d.Content
$tail.Find.Execute(" {client}", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$tail.Bold = 1
$tail.Bold = 0

# --- Change 3 -------------------------------------------------------
# Searching across the three runs " (" + "{TVA}" + "), as retainer ..."
# and replacing with the same combined text merges them into one run
# (they already share identical rPr, so no split occurs).
$d.Content.Find.Execute(" ({TVA}), as retainer regarding your case. ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " ({TVA}), as retainer regarding your case. ", 2) | Out-Null
